# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Reorder country names (Kuwait moves up, before China; Croacia moves up, before Grecia) ---
$ws.Range("A39").Value = "Kuwait"
$ws.Range("A40").Value = "China"
$ws.Range("A41").Value = "Belgica"

$ws.Range("A91").Value = "Croacia"
$ws.Range("A92").Value = "Grecia"

# --- Updated timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 31 de Agosto de 2020 a las 12:18"

# --- Updated numeric data (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6175008
$ws.Range("C4").Value = 1772
$ws.Range("E4").Value = 2561967
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 187227

# Row 15 - Iran
$ws.Range("B15").Value = 375212
$ws.Range("C15").Value = 1642
$ws.Range("D15").Value = 323233
$ws.Range("E15").Value = 30408
$ws.Range("G15").Value = 109
$ws.Range("H15").Value = 21571

# Row 18 - Banglades
$ws.Range("B18").Value = 312996
$ws.Range("C18").Value = 2174
$ws.Range("D18").Value = 204887
$ws.Range("E18").Value = 103828
$ws.Range("G18").Value = 33
$ws.Range("H18").Value = 4281

# Row 37 - Rumania
$ws.Range("B37").Value = 87540
$ws.Range("C37").Value = 755
$ws.Range("D37").Value = 37869
$ws.Range("E37").Value = 46050
$ws.Range("G37").Value = 43
$ws.Range("H37").Value = 3621

# Row 38 - Oman
$ws.Range("B38").Value = 85722
$ws.Range("C38").Value = 178
$ws.Range("D38").Value = 80810
$ws.Range("E38").Value = 4227
$ws.Range("G38").Value = 8
$ws.Range("H38").Value = 685

# Row 39 - Kuwait (new position)
$ws.Range("B39").Value = 85109
$ws.Range("C39").Value = 473
$ws.Range("D39").Value = 77224
$ws.Range("E39").Value = 7354
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = 531

# Row 40 - China (new position)
$ws.Range("B40").Value = 85048
$ws.Range("C40").Value = 17
$ws.Range("D40").Value = 80177
$ws.Range("E40").Value = 237
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 4634

# Row 41 - Belgica (new position)
$ws.Range("B41").Value = 85042
$ws.Range("C41").Value = 443
$ws.Range("D41").Value = 18415
$ws.Range("E41").Value = 56733
$ws.Range("G41").Value = 3
$ws.Range("H41").Value = 9894

# Row 71 - Austria
$ws.Range("B71").Value = 27438
$ws.Range("C71").Value = 272
$ws.Range("D71").Value = 23226
$ws.Range("E71").Value = 3479

# Row 91 - Croacia (new position)
$ws.Range("B91").Value = 10269
$ws.Range("C91").Value = 146
$ws.Range("D91").Value = 7434
$ws.Range("E91").Value = 2649
$ws.Range("G91").Value = 2
$ws.Range("H91").Value = 186

# Row 92 - Grecia (new position)
$ws.Range("B92").Value = 10134
$ws.Range("D92").Value = 3804
$ws.Range("E92").Value = 6068
$ws.Range("H92").Value = 262

# Row 100 - Haiti
$ws.Range("B100").Value = 8224
$ws.Range("C100").Value = 15
$ws.Range("E100").Value = 2195

# Row 127 - Sri Lanka
$ws.Range("B127").Value = 3015
$ws.Range("C127").Value = 3
$ws.Range("D127").Value = 2868
$ws.Range("E127").Value = 135

# Row 149 - Reunion
$ws.Range("E149").Value = 745
$ws.Range("H149").Value = 9

# Row 152 - Republica de Chipre
$ws.Range("D152").Value = 1139
$ws.Range("E152").Value = 328

# Row 179 - Islas Feroe
$ws.Range("D179").Value = 374
$ws.Range("E179").Value = 37
